$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-23 17:19:01"
$ws.Range("K2").Value = "12.9 MJ/m2"
$ws.Range("E3").Value = "2026-02-23 17:19:04"
$ws.Range("K3").Value = "16.6 MJ/m2"
$ws.Range("O3").Value = "3.9 °C"
$ws.Range("E4").Value = "2026-02-23 17:19:07"
$ws.Range("K4").Value = "14.8 MJ/m2"
$ws.Range("O4").Value = "12.4 °C"
$ws.Range("E5").Value = "2026-02-23 17:19:09"
$ws.Range("H5").Value = "'31%"
$ws.Range("K5").Value = "16.3 MJ/m2"
$ws.Range("O5").Value = "4.5 °C"
$ws.Range("E6").Value = "2026-02-23 17:19:12"
$ws.Range("J6").Value = "1024.6 hPa"
$ws.Range("K6").Value = "15.5 MJ/m2"
$ws.Range("O6").Value = "14.0 °C"
$ws.Range("E7").Value = "2026-02-23 17:19:14"
$ws.Range("J7").Value = "1024.9 hPa"
$ws.Range("K7").Value = "15.7 MJ/m2"
$ws.Range("E8").Value = "2026-02-23 17:19:16"
$ws.Range("H8").Value = "'56%"
$ws.Range("K8").Value = "15.8 MJ/m2"
$ws.Range("E9").Value = "2026-02-23 17:19:19"
$ws.Range("K9").Value = "15.2 MJ/m2"
$ws.Range("O9").Value = "12.7 °C"
$ws.Range("E10").Value = "2026-02-23 17:19:22"
$ws.Range("K10").Value = "15.4 MJ/m2"
$ws.Range("O10").Value = "11.2 °C"
$ws.Range("E11").Value = "2026-02-23 17:19:24"
$ws.Range("H11").Value = "'69%"
$ws.Range("O11").Value = "8.9 °C"
$ws.Range("E12").Value = "2026-02-23 17:19:27"
$ws.Range("O12").Value = "10.8 °C"
$ws.Range("E13").Value = "2026-02-23 17:19:30"
$ws.Range("H13").Value = "'60%"
$ws.Range("J13").Value = "1027.3 hPa"
$ws.Range("K13").Value = "15.7 MJ/m2"
$ws.Range("O13").Value = "6.7 °C"
$ws.Range("E14").Value = "2026-02-23 17:19:33"
$ws.Range("K14").Value = "15.0 MJ/m2"
$ws.Range("O14").Value = "13.2 °C"
$ws.Range("E15").Value = "2026-02-23 17:19:35"
$ws.Range("O15").Value = "12.9 °C"
$ws.Range("E16").Value = "2026-02-23 17:19:37"
$ws.Range("H16").Value = "'18%"
$ws.Range("O16").Value = "3.7 °C"
$ws.Range("E17").Value = "2026-02-23 17:19:40"
$ws.Range("K17").Value = "17.2 MJ/m2"
$ws.Range("E18").Value = "2026-02-23 17:19:43"
$ws.Range("H18").Value = "'73%"
$ws.Range("J18").Value = "1025.1 hPa"
$ws.Range("K18").Value = "15.7 MJ/m2"
$ws.Range("O18").Value = "10.8 °C"
$ws.Range("E19").Value = "2026-02-23 17:19:46"
$ws.Range("K19").Value = "15.3 MJ/m2"
$ws.Range("O19").Value = "12.6 °C"
$ws.Range("E20").Value = "2026-02-23 17:19:49"
$ws.Range("K20").Value = "16.6 MJ/m2"
$ws.Range("E21").Value = "2026-02-23 17:19:51"
$ws.Range("J21").Value = "1026.1 hPa"
$ws.Range("K21").Value = "16.3 MJ/m2"
$ws.Range("O21").Value = "9.3 °C"
$ws.Range("E22").Value = "2026-02-23 17:19:54"
$ws.Range("H22").Value = "'20%"
$ws.Range("K22").Value = "16.9 MJ/m2"
$ws.Range("E23").Value = "2026-02-23 17:19:57"
$ws.Range("K23").Value = "16.3 MJ/m2"
$ws.Range("E24").Value = "2026-02-23 17:19:59"
$ws.Range("H24").Value = "'83%"
$ws.Range("J24").Value = "1026.5 hPa"
$ws.Range("K24").Value = "16.1 MJ/m2"
$ws.Range("M24").Value = "18.3 °C 16:43 TU"
$ws.Range("O24").Value = "8.3 °C"
$ws.Range("E25").Value = "2026-02-23 17:20:02"
$ws.Range("K25").Value = "17.1 MJ/m2"
$ws.Range("O25").Value = "6.2 °C"
$ws.Range("E26").Value = "2026-02-23 17:20:05"
$ws.Range("K26").Value = "15.8 MJ/m2"
$ws.Range("E27").Value = "2026-02-23 17:20:08"
$ws.Range("K27").Value = "16.9 MJ/m2"
$ws.Range("E28").Value = "2026-02-23 17:20:10"
$ws.Range("J28").Value = "1025.1 hPa"
$ws.Range("K28").Value = "15.0 MJ/m2"
$ws.Range("O28").Value = "11.1 °C"
$ws.Range("E29").Value = "2026-02-23 17:20:13"
$ws.Range("K29").Value = "15.6 MJ/m2"
$ws.Range("O29").Value = "10.8 °C"
$ws.Range("E30").Value = "2026-02-23 17:20:16"
$ws.Range("K30").Value = "15.3 MJ/m2"
$ws.Range("E31").Value = "2026-02-23 17:20:18"
$ws.Range("J31").Value = "1024.0 hPa"
$ws.Range("K31").Value = "15.2 MJ/m2"
$ws.Range("E32").Value = "2026-02-23 17:20:21"
$ws.Range("H32").Value = "'65%"
$ws.Range("K32").Value = "15.9 MJ/m2"
$ws.Range("O32").Value = "8.4 °C"
$ws.Range("E33").Value = "2026-02-23 17:20:24"
$ws.Range("H33").Value = "'46%"
$ws.Range("J33").Value = "1025.6 hPa"
$ws.Range("K33").Value = "16.1 MJ/m2"
$ws.Range("O33").Value = "8.5 °C"
$ws.Range("E34").Value = "2026-02-23 17:20:27"
$ws.Range("H34").Value = "'37%"
$ws.Range("E35").Value = "2026-02-23 17:20:29"
$ws.Range("J35").Value = "1025.0 hPa"
$ws.Range("K35").Value = "16.8 MJ/m2"
$ws.Range("O35").Value = "12.9 °C"
$ws.Range("E36").Value = "2026-02-23 17:20:32"
$ws.Range("H36").Value = "'73%"
$ws.Range("J36").Value = "1024.8 hPa"
$ws.Range("K36").Value = "15.2 MJ/m2"
$ws.Range("O36").Value = "12.9 °C"
$ws.Range("E37").Value = "2026-02-23 17:20:35"
$ws.Range("J37").Value = "1026.7 hPa"
$ws.Range("O37").Value = "9.5 °C"
$ws.Range("E38").Value = "2026-02-23 17:20:38"
$ws.Range("K38").Value = "15.9 MJ/m2"
$ws.Range("O38").Value = "12.3 °C"
$ws.Range("E39").Value = "2026-02-23 17:20:40"
$ws.Range("E40").Value = "2026-02-23 17:20:43"
$ws.Range("H40").Value = "'61%"
$ws.Range("J40").Value = "1026.4 hPa"
$ws.Range("O40").Value = "8.9 °C"
$ws.Range("E41").Value = "2026-02-23 17:20:45"
$ws.Range("K41").Value = "15.8 MJ/m2"
$ws.Range("O41").Value = "12.4 °C"
$ws.Range("E42").Value = "2026-02-23 17:20:48"
$ws.Range("H42").Value = "'77%"
$ws.Range("O42").Value = "11.8 °C"
$ws.Range("E43").Value = "2026-02-23 17:20:50"
$ws.Range("K43").Value = "15.4 MJ/m2"
$ws.Range("O43").Value = "10.0 °C"
$ws.Range("E44").Value = "2026-02-23 17:20:53"
$ws.Range("K44").Value = "16.2 MJ/m2"
$ws.Range("E45").Value = "2026-02-23 17:20:56"
$ws.Range("J45").Value = "1026.9 hPa"
$ws.Range("E46").Value = "2026-02-23 17:20:59"
$ws.Range("H46").Value = "'74%"
$ws.Range("J46").Value = "1026.4 hPa"
$ws.Range("K46").Value = "15.6 MJ/m2"
$ws.Range("O46").Value = "9.8 °C"
